$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Value = "x"
function RGB($r,$g,$b) { return $r + ($g*256) + ($b*65536) }
$ws.Range("Z1").Interior.Color = (RGB 255 0 0)
